# Refresh the cryptocurrency price / 1h-volume-change snapshot.
#
# The Price column (D) stores numeric-looking text (e.g. "0.995",
# "50.862.10") as plain text, exactly as scraped from the source site.
# A leading apostrophe is used (the normal Excel "treat as text" marker)
# so values that look like numbers are not silently converted to the
# Number type when assigned through .Value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'50.862.10"
$ws.Range("E2").Value = "  -1.22%  "

$ws.Range("D3").Value = "'2.906.78"
$ws.Range("E3").Value = "  -0.89%  "

$ws.Range("D4").Value = "'0.995"
$ws.Range("E4").Value = "  -0.53%  "

$ws.Range("D5").Value = "'368.12"
$ws.Range("E5").Value = "  +5.12%  "

$ws.Range("D6").Value = "'102.43"
$ws.Range("E6").Value = "  -3.70%  "

$ws.Range("D7").Value = "'0.539"
$ws.Range("E7").Value = "  -2.62%  "

$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = "  -0.25%  "

$ws.Range("D9").Value = "'0.581"
$ws.Range("E9").Value = "  -3.46%  "

$ws.Range("D10").Value = "'36.67"
$ws.Range("E10").Value = "  -2.75%  "

$ws.Range("E11").Value = "  +0.75%  "

$ws.Range("D12").Value = "'0.0831"
$ws.Range("E12").Value = "  -1.85%  "

$ws.Range("D13").Value = "'18.24"
$ws.Range("E13").Value = "  -3.58%  "

$ws.Range("D14").Value = "'3.357.93"
$ws.Range("E14").Value = "  -1.11%  "

$ws.Range("D15").Value = "'7.35"
$ws.Range("E15").Value = "  -2.08%  "

$ws.Range("D16").Value = "'2.903.67"
$ws.Range("E16").Value = "  -0.93%  "

$ws.Range("D17").Value = "'0.922"
$ws.Range("E17").Value = "  -4.13%  "

$ws.Range("D18").Value = "'50.552.13"
$ws.Range("E18").Value = "  -1.79%  "

$ws.Range("D19").Value = "'3.20"
$ws.Range("E19").Value = "  -5.31%  "

$ws.Range("D20").Value = "'7.16"
$ws.Range("E20").Value = "  -2.91%  "

$ws.Range("D21").Value = "'12.86"
$ws.Range("E21").Value = "  -3.83%  "

$ws.Range("D22").Value = "'0.0₃0939"
$ws.Range("E22").Value = "  -2.44%  "

$ws.Range("D23").Value = "'67.94"
$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("D24").Value = "'257.73"
$ws.Range("E24").Value = "  -0.96%  "

$ws.Range("D25").Value = "'2.67"
$ws.Range("E25").Value = "  -0.96%  "

$ws.Range("D26").Value = "'4.21"
$ws.Range("E26").Value = "  -1.52%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("E28").Value = "  -4.21%  "

$ws.Range("D29").Value = "'25.51"
$ws.Range("E29").Value = "  -3.37%  "

$ws.Range("D30").Value = "'7.06"
$ws.Range("E30").Value = "  -2.89%  "

$ws.Range("E31").Value = "  -4.94%  "

$ws.Range("D32").Value = "'6.26"
$ws.Range("E32").Value = "  +3.32%  "

$ws.Range("D33").Value = "'9.84"
$ws.Range("E33").Value = "  -3.63%  "

$ws.Range("D35").Value = "'51.12"
$ws.Range("E35").Value = "  +1.28%  "

$ws.Range("D36").Value = "'34.17"
$ws.Range("E36").Value = "  -3.69%  "

$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("D38").Value = "'0.0420"
$ws.Range("E38").Value = "  -2.03%  "

$ws.Range("D39").Value = "'2.97"
$ws.Range("E39").Value = "  -4.81%  "

$ws.Range("D40").Value = "'16.99"
$ws.Range("E40").Value = "  -3.56%  "

$ws.Range("D41").Value = "'2.59"
$ws.Range("E41").Value = "  -1.95%  "

$ws.Range("E42").Value = "  -5.06%  "

$ws.Range("D43").Value = "'0.112"
$ws.Range("E43").Value = "  -2.80%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'21.96"
$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'118.96"
$ws.Range("E45").Value = "  -1.23%  "

$ws.Range("E46").Value = "  -2.10%  "

$ws.Range("E47").Value = "  +0.83%  "

$ws.Range("D48").Value = "'2.014.77"
$ws.Range("E48").Value = "  -3.93%  "

$ws.Range("D49").Value = "'3.13"
$ws.Range("E49").Value = "  -5.05%  "

$ws.Range("D50").Value = "'3.185.58"
$ws.Range("E50").Value = "  -0.99%  "

$ws.Range("D51").Value = "'0.236"
$ws.Range("E51").Value = "  -0.19%  "
